$d = $word.ActiveDocument

$pairs = @(
    @("2024-07-18 Thursday", "2024-07-19 Friday"),
    @("834×3=2502", "362×5=1810"),
    @("248×9=2232", "392×4=1568"),
    @("689×8=5512", "792×9=7128"),
    @("131×8=1048", "262×6=1572"),
    @("785×7=5495", "409×7=2863"),
    @("561×8=4488", "264×3=792"),
    @("194×4=776", "345×9=3105"),
    @("706×3=2118", "957×8=7656"),
    @("751×4=3004", "184×2=368"),
    @("548×7=3836", "698×5=3490"),
    @("956×6=5736", "960×8=7680"),
    @("968×2=1936", "736×2=1472"),
    @("680×4=2720", "992×7=6944"),
    @("807×7=5649", "305×5=1525"),
    @("509×4=2036", "302×3=906"),
    @("222×5=1110", "827×2=1654"),
    @("803×6=4818", "648×6=3888"),
    @("318×3=954", "858×2=1716"),
    @("312×8=2496", "917×9=8253"),
    @("653×2=1306", "941×8=7528"),
    @("357×6=2142", "651×7=4557"),
    @("838×3=2514", "280×9=2520"),
    @("544×2=1088", "654×9=5886"),
    @("954×6=5724", "524×8=4192"),
    @("679×4=2716", "272×9=2448")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $found = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND: $old"
    }
}

Write-Host "Done"
